$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.495728850364685
$ws.Range("B1").Value = 1.793264746665955
$ws.Range("C1").Value = 2.403722524642944
$ws.Range("D1").Value = 5.025171279907227
$ws.Range("E1").Value = 1.453248620033264
